# Update FFXIV Leve profit calculations (currentAveragePrice / LevePriceNQ|HQ / LeveProfitNQ|HQ)
# across the 8 crafting-class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below are data refreshes (market price re-pull); no formulas are involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 5968.684   # H9: 6295.8335 -> 5968.684
$ws.Cells.Item(9, 10).Value = 17353.334   # J9: 20808 -> 17353.334
$ws.Cells.Item(9, 12).Value = 17353.334   # L9: 20808 -> 17353.334
$ws.Cells.Item(9, 14).Value = -17691.334   # N9: -21146 -> -17691.334
$ws.Cells.Item(40, 8).Value = 5371.143   # H40: 4912.125 -> 5371.143
$ws.Cells.Item(40, 9).Value = 4679.6   # I40: 4182.8335 -> 4679.6
$ws.Cells.Item(40, 11).Value = 4679.6   # K40: 4182.8335 -> 4679.6
$ws.Cells.Item(40, 13).Value = -4504.6   # M40: -4007.8335 -> -4504.6
$ws.Cells.Item(55, 8).Value = 271.85715   # H55: 418.5909 -> 271.85715
$ws.Cells.Item(55, 9).Value = 129.83333   # I55: 389.07693 -> 129.83333
$ws.Cells.Item(55, 11).Value = 129.83333   # K55: 389.07693 -> 129.83333
$ws.Cells.Item(55, 13).Value = 84.16667000000001   # M55: -175.07693 -> 84.16667000000001
$ws.Cells.Item(59, 8).Value = 3393.4   # H59: 2995 -> 3393.4
$ws.Cells.Item(59, 10).Value = 4989   # J59: 4990 -> 4989
$ws.Cells.Item(59, 12).Value = 14967   # L59: 14970 -> 14967
$ws.Cells.Item(59, 14).Value = -16081   # N59: -16084 -> -16081
$ws.Cells.Item(116, 8).Value = 7470.2   # H116: 7603 -> 7470.2
$ws.Cells.Item(116, 9).Value = 6948.5713   # I116: 7116.1665 -> 6948.5713
$ws.Cells.Item(116, 11).Value = 6948.5713   # K116: 7116.1665 -> 6948.5713
$ws.Cells.Item(116, 13).Value = -3506.5713   # M116: -3674.1665 -> -3506.5713
$ws.Cells.Item(125, 8).Value = 2426.7144   # H125: 2499.8 -> 2426.7144
$ws.Cells.Item(125, 9).Value = 1990   # I125: 0 -> 1990
$ws.Cells.Item(125, 10).Value = 2499.5   # J125: 2499.8 -> 2499.5
$ws.Cells.Item(125, 11).Value = 17910   # K125: 0 -> 17910
$ws.Cells.Item(125, 12).Value = 22495.5   # L125: 22498.2 -> 22495.5
$ws.Cells.Item(125, 13).Value = -15450   # M125: (empty) -> -15450
$ws.Cells.Item(125, 14).Value = -27415.5   # N125: -27418.2 -> -27415.5
$ws.Cells.Item(132, 8).Value = 1936.7609   # H132: 1972.7954 -> 1936.7609
$ws.Cells.Item(132, 10).Value = 1144   # J132: 0 -> 1144
$ws.Cells.Item(132, 12).Value = 3432   # L132: 0 -> 3432
$ws.Cells.Item(132, 14).Value = -8492   # N132: (empty) -> -8492
$ws.Cells.Item(137, 8).Value = 133133.02   # H137: 138550.69 -> 133133.02
$ws.Cells.Item(137, 9).Value = 1871.7317   # I137: 1947.2051 -> 1871.7317
$ws.Cells.Item(137, 11).Value = 5615.1951   # K137: 5841.615299999999 -> 5615.1951
$ws.Cells.Item(137, 13).Value = -3065.1951   # M137: -3291.615299999999 -> -3065.1951
$ws.Cells.Item(138, 8).Value = 2558.2046   # H138: 2561.6667 -> 2558.2046
$ws.Cells.Item(138, 9).Value = 1874.875   # I138: 1834.8529 -> 1874.875
$ws.Cells.Item(138, 10).Value = 4380.4165   # J138: 4326.7856 -> 4380.4165
$ws.Cells.Item(138, 11).Value = 5624.625   # K138: 5504.5587 -> 5624.625
$ws.Cells.Item(138, 12).Value = 13141.2495   # L138: 12980.3568 -> 13141.2495
$ws.Cells.Item(138, 13).Value = -484.625   # M138: -364.5587000000005 -> -484.625
$ws.Cells.Item(138, 14).Value = -23421.2495   # N138: -23260.3568 -> -23421.2495
$ws.Cells.Item(141, 8).Value = 1743.7858   # H141: 1698.8276 -> 1743.7858
$ws.Cells.Item(141, 9).Value = 591.04346   # I141: 584.75 -> 591.04346
$ws.Cells.Item(141, 11).Value = 1773.13038   # K141: 1754.25 -> 1773.13038
$ws.Cells.Item(141, 13).Value = 3406.86962   # M141: 3425.75 -> 3406.86962

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4277.4067   # H32: 4194.527 -> 4277.4067
$ws.Cells.Item(32, 9).Value = 2540.6626   # I32: 2490.8472 -> 2540.6626
$ws.Cells.Item(32, 11).Value = 2540.6626   # K32: 2490.8472 -> 2540.6626
$ws.Cells.Item(32, 13).Value = -2253.6626   # M32: -2203.8472 -> -2253.6626
$ws.Cells.Item(132, 8).Value = 8152.1523   # H132: 8317.111000000001 -> 8152.1523
$ws.Cells.Item(132, 9).Value = 9326.757   # I132: 9565.583000000001 -> 9326.757
$ws.Cells.Item(132, 11).Value = 27980.271   # K132: 28696.749 -> 27980.271
$ws.Cells.Item(132, 13).Value = -25450.271   # M132: -26166.749 -> -25450.271

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1328.129   # H20: 1331.3549 -> 1328.129
$ws.Cells.Item(20, 9).Value = 1343.5834   # I20: 1397.6522 -> 1343.5834
$ws.Cells.Item(20, 10).Value = 1275.1428   # J20: 1140.75 -> 1275.1428
$ws.Cells.Item(20, 11).Value = 1343.5834   # K20: 1397.6522 -> 1343.5834
$ws.Cells.Item(20, 12).Value = 1275.1428   # L20: 1140.75 -> 1275.1428
$ws.Cells.Item(20, 13).Value = -1096.5834   # M20: -1150.6522 -> -1096.5834
$ws.Cells.Item(20, 14).Value = -1769.1428   # N20: -1634.75 -> -1769.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 13951.637   # H15: 15100.9 -> 13951.637
$ws.Cells.Item(15, 10).Value = 13951.637   # J15: 15100.9 -> 13951.637
$ws.Cells.Item(15, 12).Value = 13951.637   # L15: 15100.9 -> 13951.637
$ws.Cells.Item(15, 14).Value = -14291.637   # N15: -15440.9 -> -14291.637
$ws.Cells.Item(22, 8).Value = 223.83333   # H22: 206.6842 -> 223.83333
$ws.Cells.Item(22, 9).Value = 159.92857   # I22: 142.46666 -> 159.92857
$ws.Cells.Item(22, 11).Value = 159.92857   # K22: 142.46666 -> 159.92857
$ws.Cells.Item(22, 13).Value = 190.07143   # M22: 207.53334 -> 190.07143
$ws.Cells.Item(31, 8).Value = 2568.7969   # H31: 2584.1746 -> 2568.7969
$ws.Cells.Item(31, 9).Value = 2191.5   # I31: 2215.16 -> 2191.5
$ws.Cells.Item(31, 11).Value = 2191.5   # K31: 2215.16 -> 2191.5
$ws.Cells.Item(31, 13).Value = -1896.5   # M31: -1920.16 -> -1896.5
$ws.Cells.Item(34, 8).Value = 2568.7969   # H34: 2584.1746 -> 2568.7969
$ws.Cells.Item(34, 9).Value = 2191.5   # I34: 2215.16 -> 2191.5
$ws.Cells.Item(34, 11).Value = 2191.5   # K34: 2215.16 -> 2191.5
$ws.Cells.Item(34, 13).Value = -1989.5   # M34: -2013.16 -> -1989.5
$ws.Cells.Item(58, 8).Value = 3697.7097   # H58: 3860.2903 -> 3697.7097
$ws.Cells.Item(58, 9).Value = 3266.4783   # I58: 3284.76 -> 3266.4783
$ws.Cells.Item(58, 10).Value = 4937.5   # J58: 6258.3335 -> 4937.5
$ws.Cells.Item(58, 11).Value = 3266.4783   # K58: 3284.76 -> 3266.4783
$ws.Cells.Item(58, 12).Value = 4937.5   # L58: 6258.3335 -> 4937.5
$ws.Cells.Item(58, 13).Value = -3063.4783   # M58: -3081.76 -> -3063.4783
$ws.Cells.Item(58, 14).Value = -5343.5   # N58: -6664.3335 -> -5343.5
$ws.Cells.Item(99, 8).Value = 4242.5   # H99: 4368.846 -> 4242.5
$ws.Cells.Item(99, 9).Value = 4338.923   # I99: 4483.8335 -> 4338.923
$ws.Cells.Item(99, 11).Value = 4338.923   # K99: 4483.8335 -> 4338.923
$ws.Cells.Item(99, 13).Value = -2840.923   # M99: -2985.8335 -> -2840.923
$ws.Cells.Item(122, 8).Value = 5657   # H122: 3278.5 -> 5657
$ws.Cells.Item(122, 9).Value = 0   # I122: 900 -> 0
$ws.Cells.Item(122, 11).Value = 0   # K122: 2700 -> 0
$ws.Cells.Item(122, 13).ClearContents()   # M122: -250 -> (removed)
$ws.Cells.Item(126, 8).Value = 4242.5   # H126: 4368.846 -> 4242.5
$ws.Cells.Item(126, 9).Value = 4338.923   # I126: 4483.8335 -> 4338.923
$ws.Cells.Item(126, 11).Value = 13016.769   # K126: 13451.5005 -> 13016.769
$ws.Cells.Item(126, 13).Value = -10546.769   # M126: -10981.5005 -> -10546.769
$ws.Cells.Item(132, 8).Value = 5456.552   # H132: 5601.4287 -> 5456.552
$ws.Cells.Item(132, 9).Value = 2448.3076   # I132: 2490.24 -> 2448.3076
$ws.Cells.Item(132, 11).Value = 7344.9228   # K132: 7470.719999999999 -> 7344.9228
$ws.Cells.Item(132, 13).Value = -4814.9228   # M132: -4940.719999999999 -> -4814.9228
$ws.Cells.Item(136, 8).Value = 3697.7097   # H136: 3860.2903 -> 3697.7097
$ws.Cells.Item(136, 9).Value = 3266.4783   # I136: 3284.76 -> 3266.4783
$ws.Cells.Item(136, 10).Value = 4937.5   # J136: 6258.3335 -> 4937.5
$ws.Cells.Item(136, 11).Value = 9799.4349   # K136: 9854.280000000001 -> 9799.4349
$ws.Cells.Item(136, 12).Value = 14812.5   # L136: 18775.0005 -> 14812.5
$ws.Cells.Item(136, 13).Value = -7249.4349   # M136: -7304.280000000001 -> -7249.4349
$ws.Cells.Item(136, 14).Value = -19912.5   # N136: -23875.0005 -> -19912.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 499.5   # H8: 500 -> 499.5
$ws.Cells.Item(8, 9).Value = 499.5   # I8: 500 -> 499.5
$ws.Cells.Item(8, 11).Value = 1498.5   # K8: 1500 -> 1498.5
$ws.Cells.Item(8, 13).Value = -1359.5   # M8: -1361 -> -1359.5
$ws.Cells.Item(9, 8).Value = 145185   # H9: 185198.33 -> 145185
$ws.Cells.Item(9, 9).Value = 191865   # I9: 275225 -> 191865
$ws.Cells.Item(9, 11).Value = 575595   # K9: 825675 -> 575595
$ws.Cells.Item(9, 13).Value = -575371   # M9: -825451 -> -575371
$ws.Cells.Item(131, 8).Value = 3401.8518   # H131: 2879 -> 3401.8518
$ws.Cells.Item(131, 9).Value = 2780   # I131: 1469.8182 -> 2780
$ws.Cells.Item(131, 10).Value = 3543.182   # J131: 3617.1428 -> 3543.182
$ws.Cells.Item(131, 11).Value = 8340   # K131: 4409.4546 -> 8340
$ws.Cells.Item(131, 12).Value = 10629.546   # L131: 10851.4284 -> 10629.546
$ws.Cells.Item(131, 13).Value = -3300   # M131: 630.5454 -> -3300
$ws.Cells.Item(131, 14).Value = -20709.546   # N131: -20931.4284 -> -20709.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 102288.445   # H122: 92229.39999999999 -> 102288.445
$ws.Cells.Item(122, 9).Value = 130515.14   # I122: 114413 -> 130515.14
$ws.Cells.Item(122, 11).Value = 391545.42   # K122: 343239 -> 391545.42
$ws.Cells.Item(122, 13).Value = -389095.42   # M122: -340789 -> -389095.42
$ws.Cells.Item(132, 8).Value = 6253.7144   # H132: 7048.8335 -> 6253.7144
$ws.Cells.Item(132, 9).Value = 6253.7144   # I132: 7048.8335 -> 6253.7144
$ws.Cells.Item(132, 11).Value = 18761.1432   # K132: 21146.5005 -> 18761.1432
$ws.Cells.Item(132, 13).Value = -16231.1432   # M132: -18616.5005 -> -16231.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3772.111   # H7: 3929.625 -> 3772.111
$ws.Cells.Item(7, 9).Value = 3573   # I7: 3749.8333 -> 3573
$ws.Cells.Item(7, 11).Value = 3573   # K7: 3749.8333 -> 3573
$ws.Cells.Item(7, 13).Value = -3461   # M7: -3637.8333 -> -3461
$ws.Cells.Item(30, 8).Value = 23168.7   # H30: 13922.667 -> 23168.7
$ws.Cells.Item(30, 9).Value = 36649.332   # I30: 6688.3335 -> 36649.332
$ws.Cells.Item(30, 10).Value = 17391.285   # J30: 17539.834 -> 17391.285
$ws.Cells.Item(30, 11).Value = 36649.332   # K30: 6688.3335 -> 36649.332
$ws.Cells.Item(30, 12).Value = 17391.285   # L30: 17539.834 -> 17391.285
$ws.Cells.Item(30, 13).Value = -36541.332   # M30: -6580.3335 -> -36541.332
$ws.Cells.Item(30, 14).Value = -17607.285   # N30: -17755.834 -> -17607.285
$ws.Cells.Item(100, 8).Value = 3207.4   # H100: 3282.88 -> 3207.4
$ws.Cells.Item(100, 9).Value = 2468.682   # I100: 2469.1365 -> 2468.682
$ws.Cells.Item(100, 10).Value = 8624.666999999999   # J100: 9250.333000000001 -> 8624.666999999999
$ws.Cells.Item(100, 11).Value = 2468.682   # K100: 2469.1365 -> 2468.682
$ws.Cells.Item(100, 12).Value = 8624.666999999999   # L100: 9250.333000000001 -> 8624.666999999999
$ws.Cells.Item(100, 13).Value = -1927.682   # M100: -1928.1365 -> -1927.682
$ws.Cells.Item(100, 14).Value = -9706.666999999999   # N100: -10332.333 -> -9706.666999999999
$ws.Cells.Item(126, 8).Value = 3772.111   # H126: 3929.625 -> 3772.111
$ws.Cells.Item(126, 9).Value = 3573   # I126: 3749.8333 -> 3573
$ws.Cells.Item(126, 11).Value = 10719   # K126: 11249.4999 -> 10719
$ws.Cells.Item(126, 13).Value = -8249   # M126: -8779.499899999999 -> -8249
$ws.Cells.Item(132, 8).Value = 2624.7163   # H132: 2591.2236 -> 2624.7163
$ws.Cells.Item(132, 9).Value = 2155.2036   # I132: 2126.5178 -> 2155.2036
$ws.Cells.Item(132, 11).Value = 6465.610799999999   # K132: 6379.553400000001 -> 6465.610799999999
$ws.Cells.Item(132, 13).Value = -3935.610799999999   # M132: -3849.553400000001 -> -3935.610799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3353.7346   # H122: 3270.157 -> 3353.7346
$ws.Cells.Item(122, 9).Value = 3217.3684   # I122: 3117.625 -> 3217.3684
$ws.Cells.Item(122, 11).Value = 9652.1052   # K122: 9352.875 -> 9652.1052
$ws.Cells.Item(122, 13).Value = -7202.1052   # M122: -6902.875 -> -7202.1052
$ws.Cells.Item(126, 8).Value = 1972.9688   # H126: 1973 -> 1972.9688
$ws.Cells.Item(126, 9).Value = 1975.8966   # I126: 1975.931 -> 1975.8966
$ws.Cells.Item(126, 11).Value = 5927.6898   # K126: 5927.793 -> 5927.6898
$ws.Cells.Item(126, 13).Value = -3457.6898   # M126: -3457.793 -> -3457.6898
$ws.Cells.Item(136, 8).Value = 4823.518   # H136: 4837.982 -> 4823.518
$ws.Cells.Item(136, 9).Value = 5103.8   # I136: 5124.05 -> 5103.8
$ws.Cells.Item(136, 11).Value = 15311.4   # K136: 15372.15 -> 15311.4
$ws.Cells.Item(136, 13).Value = -12761.4   # M136: -12822.15 -> -12761.4
